# Update the "dSF" column (F) values for a set of rows.
# These values were repulled / recalculated (mean calculation) per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    4  = -5
    6  = 5
    8  = 5
    10 = 6
    14 = -4
    17 = -1
    19 = 2
    20 = 1
    24 = 1
    26 = 0
    28 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
